$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 3 (Ontario) figures ---
$ws.Range("G3").Value = 25470.31
$ws.Range("H3").Value = 27139.01
$ws.Range("I3").Value = -3871.81

# --- Append a new row 4 (same shape/values as row 2, but ID = 3) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 30000
$ws.Range("C4").Value = "alb"
$ws.Range("D4").Value = 32000
$ws.Range("E4").Value = 23267.2
$ws.Range("F4").Value = 32550
$ws.Range("G4").Value = 23667.1
$ws.Range("H4").Value = 25923.77
$ws.Range("I4").Value = -2656.57

# Row 4 should carry the same formatting as row 2: B/C/D left-aligned
# (style index 1), the rest left at the workbook default style.
$ws.Range("B4:D4").HorizontalAlignment = -4131
$ws.Range("A4").Style = "Normal"
$ws.Range("E4:I4").Style = "Normal"

# --- Move the active selection to E5, matching the post-edit sheet view ---
$ws.Range("E5").Select()
